$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin snapshot (price + 1h change); rows 27/28 and 45/46 swap rank order.
# "Price" values that look like plain numbers (single decimal point, e.g. "19.00")
# are written with a temporary Text format so Excel keeps the exact digits/trailing
# zeros instead of silently re-parsing them as floating point numbers; the format is
# then restored to the sheet default so no visible style change is left behind.
$ws.Range("D2").Value = "56.799.26"
$ws.Range("E2").Value = "  +9.94%  "
$ws.Range("D3").Value = "3.247.37"
$ws.Range("E3").Value = "  +4.72%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "394.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.37%  "
$ws.Range("D7").Value = "3.241.37"
$ws.Range("E7").Value = "  +4.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.10%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.622"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0958"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.75%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "3.762.37"
$ws.Range("E14").Value = "  +5.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("D17").Value = "3.258.46"
$ws.Range("E17").Value = "  +5.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "56.702.25"
$ws.Range("E20").Value = "  +9.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("E22").Value = "  +8.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "299.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.70%  "
$ws.Range("E26").Value = "  -3.30%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.70%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.169"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0480"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("E37").Value = "  +2.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "133.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.119"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.35%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.17%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("E47").Value = "  -3.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "2.138.06"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("E51").Value = "  -5.43%  "
